$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for Q4-Q7 (rows 6 to 9)
$ws.Range("A6:G9").EntireRow.Delete()

# Update values for Q0 (row2)
$ws.Range("B2").Value = 0.1673564803967629
$ws.Range("C2").Value = 0.6573426014927735
$ws.Range("D2").Value = 0.9595486229626961
$ws.Range("E2").Value = 0.9795655276512624
$ws.Range("F2").Value = 1.001597417975854

# Update values for Q1 (row3)
$ws.Range("B3").Value = -0.09717246228638296
$ws.Range("C3").Value = 0.5709473934814091
$ws.Range("D3").Value = 0.6131476725437131
$ws.Range("E3").Value = 0.7830374656066675
$ws.Range("F3").Value = 0.8190137599691435
$ws.Range("G3").Value = 10

# Update values for Q2 (row4)
$ws.Range("B4").Value = -0.1587865486592277
$ws.Range("C4").Value = 0.292352465870846
$ws.Range("D4").Value = 0.1239581044362055
$ws.Range("E4").Value = 0.3520768445044427
$ws.Range("F4").Value = 0.3442294637030879
$ws.Range("G4").Value = 6

# Update values for Q3 (row5)
$ws.Range("B5").Value = 0.04921794047849792
$ws.Range("C5").Value = 0.2316309012071306
$ws.Range("D5").Value = 0.05607528005897246
$ws.Range("E5").Value = 0.2368021960602825
$ws.Range("F5").Value = 0.3275755619518266
$ws.Range("G5").Value = 2
